$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - update values
$ws.Range("A1").Value2 = 6
$ws.Range("B1").Value2 = 8
$ws.Range("C1").Value2 = 9
$ws.Range("D1").Value2 = 58
$ws.Range("E1").Value2 = 67
$ws.Range("F1").Value2 = 15
$ws.Range("G1").Value2 = 4
$ws.Range("H1").Value2 = 83
$ws.Range("I1").Value2 = 34
$ws.Range("J1").Value2 = 106
$ws.Range("K1").Value2 = 3
$ws.Range("L1").Value2 = 1

# M1 no longer used - clear it
$ws.Range("M1").ClearContents()

# New row 7
$ws.Range("A7").Value2 = 124

Write-Host "done"
